$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (NCTId), shifting existing columns right.
$ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Cells.Item(1, 2).Value = "status_label"

# Map the emoji status in column A to a French text label in the new column B
$redSquare = "🟥"
$orangeSquare = "🟧"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $status = [string]$ws.Cells.Item($r, 1).Value2
    if ($status -eq $redSquare) {
        $label = "rouge"
    } elseif ($status -eq $orangeSquare) {
        $label = "orange"
    } else {
        $label = $null
    }
    $ws.Cells.Item($r, 2).Value = $label
}
